$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update input values (C/D/E columns, rows 3-6) that drive the MOD() formulas ---
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 3

$ws.Range("D4").Value = 0

$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 2

$ws.Range("C6").Value = 3
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 3

# --- Unhide helper columns F:J and give them a real width ---
$ws.Range("F1:J1").EntireColumn.Hidden = $false
$ws.Range("F1:J1").ColumnWidth = 10.7

# --- Move the active selection from E6 to G3 ---
$ws.Range("G3").Select()
